# The underlying deck currently ships two DrawingML theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" palette (used by the Notes Master)
#   ppt/theme/theme2.xml -> "Integral"     palette (used by the Slide Master /
#                                                    the actual design applied
#                                                    to every slide)
#
# The target revision swaps those two palettes: the deck's live design
# (Slide Master / theme2.xml) becomes the standard "Office Theme" palette,
# while the old "Integral" palette is kept around (theme1.xml). The font
# scheme and format scheme (fills/lines/effects) are identical between the
# two themes already, so only the 12 theme colors actually need to change.
#
# Re-color the presentation's live theme (reached through
# SlideMaster.Theme.ThemeColorScheme) to the "Office Theme" color values, in
# MsoThemeColorSchemeIndex order:
#   1 Dark1, 2 Light1, 3 Dark2, 4 Light2,
#   5-10 Accent1..Accent6, 11 Hyperlink, 12 FollowedHyperlink

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

function Set-ThemeColor {
    param(
        [int]$Index,
        [string]$Hex
    )
    $r = [Convert]::ToInt32($Hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4, 2), 16)
    $colorScheme.Colors($Index).RGB = $r + ($g * 256) + ($b * 65536)
}

Set-ThemeColor 1  "000000"   # dk1
Set-ThemeColor 2  "FFFFFF"   # lt1
Set-ThemeColor 3  "44546A"   # dk2
Set-ThemeColor 4  "E7E6E6"   # lt2
Set-ThemeColor 5  "5B9BD5"   # accent1
Set-ThemeColor 6  "ED7D31"   # accent2
Set-ThemeColor 7  "A5A5A5"   # accent3
Set-ThemeColor 8  "FFC000"   # accent4
Set-ThemeColor 9  "4472C4"   # accent5
Set-ThemeColor 10 "70AD47"   # accent6
Set-ThemeColor 11 "0563C1"   # hlink
Set-ThemeColor 12 "954F72"   # folHlink

# Best-effort: try to relabel the design/theme to match its new "Office
# Theme" identity (a no-op on hosts where Theme.Name / Design.Name aren't
# persisted, harmless either way).
try { $master.Theme.Name = "Office Theme" } catch {}
try { $p.Designs.Item(1).Name = "Office Theme" } catch {}
